# Generate Report for Handback
# Marks the two in-flight files as handed back (in sync with en-US) on both
# the zh-cn and de-de status sheets, stamping the "Latest Handback
# DateTime" and filling in the "Latest Target File" / "Latest Handback
# File" columns (with hyperlinks) that a completed handback publishes.

$wb = $excel.ActiveWorkbook

$mdUrl22bc = "https://github.com/OpenLocalizationTest/oltest/blob/ab413cc0965d7fc75043a4c1efd1f21fd1da6b58/e2e/22bc109e-c38e-4086-909c-ae9d09168dc8.md"
$mdUrlDaab = "https://github.com/OpenLocalizationTest/oltest/blob/ab413cc0965d7fc75043a4c1efd1f21fd1da6b58/e2e/daabbfdb-628e-46ef-b14e-3a8486e2eea3.md"

$statusText = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("C2").Value = $statusText
$ws.Range("C3").Value = $statusText

$ws.Range("H2").Value = "2016-03-14 08:43:54"
$ws.Range("H3").Value = "2016-03-14 08:43:54"

$ws.Hyperlinks.Add($ws.Range("F2"), $mdUrl22bc, $null, $null, "22bc109e-c38e-4086-909c-ae9d09168dc8.md")
$ws.Hyperlinks.Add($ws.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/666501d5cea3dbe25c0f8e40969ae4c380f9e525/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/22bc109e-c38e-4086-909c-ae9d09168dc8.d04c2f87d13be762fe6b02e2d5b70f55e80abda2.zh-cn.xlf", $null, $null, "22bc109e-c38e-4086-909c-ae9d09168dc8.d04c2f87d13be762fe6b02e2d5b70f55e80abda2.zh-cn.xlf")

$ws.Hyperlinks.Add($ws.Range("F3"), $mdUrlDaab, $null, $null, "daabbfdb-628e-46ef-b14e-3a8486e2eea3.md")
$ws.Hyperlinks.Add($ws.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/666501d5cea3dbe25c0f8e40969ae4c380f9e525/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/daabbfdb-628e-46ef-b14e-3a8486e2eea3.1329f45496df605e11eaf200c985464ebf08da1b.zh-cn.xlf", $null, $null, "daabbfdb-628e-46ef-b14e-3a8486e2eea3.1329f45496df605e11eaf200c985464ebf08da1b.zh-cn.xlf")

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("C2").Value = $statusText
$ws.Range("C3").Value = $statusText

$ws.Range("H2").Value = "2016-03-14 08:44:00"
$ws.Range("H3").Value = "2016-03-14 08:44:00"

$ws.Hyperlinks.Add($ws.Range("F2"), $mdUrl22bc, $null, $null, "22bc109e-c38e-4086-909c-ae9d09168dc8.md")
$ws.Hyperlinks.Add($ws.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e100855a0eceaca414da6f53f5636ccda06e2f70/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/22bc109e-c38e-4086-909c-ae9d09168dc8.d04c2f87d13be762fe6b02e2d5b70f55e80abda2.de-de.xlf", $null, $null, "22bc109e-c38e-4086-909c-ae9d09168dc8.d04c2f87d13be762fe6b02e2d5b70f55e80abda2.de-de.xlf")

$ws.Hyperlinks.Add($ws.Range("F3"), $mdUrlDaab, $null, $null, "daabbfdb-628e-46ef-b14e-3a8486e2eea3.md")
$ws.Hyperlinks.Add($ws.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e100855a0eceaca414da6f53f5636ccda06e2f70/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/daabbfdb-628e-46ef-b14e-3a8486e2eea3.1329f45496df605e11eaf200c985464ebf08da1b.de-de.xlf", $null, $null, "daabbfdb-628e-46ef-b14e-3a8486e2eea3.1329f45496df605e11eaf200c985464ebf08da1b.de-de.xlf")
